$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-24 down to 14-25)
$ws.Rows.Item(13).Insert()

# The insert copies column A formatting from row 12 into the new blank A13; remove it (row 13 has no A cell in the target layout)
$ws.Range("A13").Clear()

# Give the new B13/C13 cells the same formatting as the other data rows (style index 2 / 3)
$ws.Range("B11").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fix up / fill in the correct cell text for the affected rows
$ws.Range("B10").Value = "Dar ao futuro engenheiro os conceitos fundamentais relacionados ao escoamento de fluidos e desenvolver as equações de conservação de massa, energia e quantidade de movimento. Os conceitos e modelos matemáticos estudados servem de base para a compreensão dos processos produtivos que envolvam a transferência de fluidos e para as disciplinas de Operações Unitárias que estudam os princípios destas operações."
$ws.Range("C10").Value = "Dar ao futuro engenheiro os conceitos fundamentais relacionados ao escoamento de fluidos e desenvolver as equações de conservação de massa, energia e quantidade de movimento. Os conceitos e modelos matemáticos estudados servem de base para a compreensão dos processos produtivos que envolvam a transferência de fluidos e para as disciplinas de Operações Unitárias que estudam os princípios destas operações."
$ws.Range("B13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("B14").Value = "1) Bases conceituais para o estudo dos Fenômenos de transporte `n2) Propriedades gerais dos fluidos `n3) Cinemática dos fluidos:. `n4) Equações de Conservação na forma Integral: `n5) Equações Diferenciais do Escoamento de Fluidos: `n6) Teoria da Camada Limite: `n7) Escoamento em tubos:"
$ws.Range("C14").Value = "1) Bases conceituais para o estudo dos Fenômenos de transporte `n2) Propriedades gerais dos fluidos `n3) Cinemática dos fluidos:. `n4) Equações de Conservação na forma Integral: `n5) Equações Diferenciais do Escoamento de Fluidos: `n6) Teoria da Camada Limite: `n7) Escoamento em tubos:"
$ws.Range("B16").Value = "1) Bases conceituais para o estudo dos Fenômenos de transporte `nFluidos e a hipótese do contínuo. Importância da análise dimensional e uso dos números adimensionais. Leis básicas para transferência de massa, calor e quantidade de movimento. Lei geral para os fenômenos de transporte. Difusividade molecular, condutividade térmica e viscosidade. Transporte simultâneo de massa, calor e quantidade de movimento. Formulação integral e diferencial. `n2) Propriedades gerais dos fluidos: Massa específica, peso específico, volume específico. Tensão e Pressão. Fluidos Newtonianos e não Newtonianos. Viscosidade. Tensão superficial e Capilaridade. Módulo de elasticidade volumétrica e compressibilidade. `n3) Cinemática dos fluidos: Descrição de um Fluido em Movimento: Método de Euler e Lagrange - Campo de escoamento de um fluido - Escoamento permanente e transiente - Trajetórias e Linhas de corrente - Sistema e volume de controle - Escoamentos unidimensionais e bidimensionais. Escoamento uniforme. Escoamento laminar e turbulento: N° de Reynolds. `n4) Equações de Conservação na forma Integral: Fluxo de uma grandeza. Conservação da Massa, continuidade. Formas específicas para a expressão integral. Conservação da quantidade de movimento linear. Conservação da Energia. Equação de Bernoulli. Aplicações `n5) Equações Diferenciais do Escoamento de Fluidos: Equação da conservação da massa e continuidade. Equação da energia. Equação de Navier-Stokes. Aplicações `n6) Teoria da Camada Limite: Definição de camada limite . Camada limite laminar e turbulenta. Camada limite sobre uma placa plana. Aplicações `n7) Escoamento em tubos: Escoamento Laminar e turbulento. Coeficiente de atrito. Região turbulenta e de transição. Diagramas de Moody e Von Karman . Equação da energia com equipamentos de transporte. Perda de carga em acidentes. Diâmetro equivalente."
$ws.Range("C16").Value = "1) Bases conceituais para o estudo dos Fenômenos de transporte `nFluidos e a hipótese do contínuo. Importância da análise dimensional e uso dos números adimensionais. Leis básicas para transferência de massa, calor e quantidade de movimento. Lei geral para os fenômenos de transporte. Difusividade molecular, condutividade térmica e viscosidade. Transporte simultâneo de massa, calor e quantidade de movimento. Formulação integral e diferencial. `n2) Propriedades gerais dos fluidos: Massa específica, peso específico, volume específico. Tensão e Pressão. Fluidos Newtonianos e não Newtonianos. Viscosidade. Tensão superficial e Capilaridade. Módulo de elasticidade volumétrica e compressibilidade. `n3) Cinemática dos fluidos: Descrição de um Fluido em Movimento: Método de Euler e Lagrange - Campo de escoamento de um fluido - Escoamento permanente e transiente - Trajetórias e Linhas de corrente - Sistema e volume de controle - Escoamentos unidimensionais e bidimensionais. Escoamento uniforme. Escoamento laminar e turbulento: N° de Reynolds. `n4) Equações de Conservação na forma Integral: Fluxo de uma grandeza. Conservação da Massa, continuidade. Formas específicas para a expressão integral. Conservação da quantidade de movimento linear. Conservação da Energia. Equação de Bernoulli. Aplicações `n5) Equações Diferenciais do Escoamento de Fluidos: Equação da conservação da massa e continuidade. Equação da energia. Equação de Navier-Stokes. Aplicações `n6) Teoria da Camada Limite: Definição de camada limite . Camada limite laminar e turbulenta. Camada limite sobre uma placa plana. Aplicações `n7) Escoamento em tubos: Escoamento Laminar e turbulento. Coeficiente de atrito. Região turbulenta e de transição. Diagramas de Moody e Von Karman . Equação da energia com equipamentos de transporte. Perda de carga em acidentes. Diâmetro equivalente."
$ws.Range("B19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("C19").Value = "Aplicação de 2 provas, P1 e P2."
$ws.Range("B20").Value = "A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)"
$ws.Range("C20").Value = "A média do período será MP = (P1+2P2)/3. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)"
$ws.Range("B21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("C21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("B22").Value = "1) YONG, D. F.; OKIISHI, T. H.; MUNSON, B.R. Fundamentos da mecânica dos fluidos. São Paulo: Edgard Blucher `n2) BRUNETTI, F. Mecânica dos fluídos. São Paulo: Pearson Education. `n3) FOX, Robert W. Introdução à mecânica dos fluídos. Rio de Janeiro: LTC. `n4) WHITE, Frank M. Mecânica dos fluídos. Rio de Janeiro: Mcgraw-hill Interamericana. `n5) COULSON, J. M.; RICHARDSON, J.F. Chemical engineering . Oxford: Butterworth Heinemann. Volume 1: Fluid Flow, Heat Transfer and Mass Transfer"
$ws.Range("C22").Value = "1) YONG, D. F.; OKIISHI, T. H.; MUNSON, B.R. Fundamentos da mecânica dos fluidos. São Paulo: Edgard Blucher `n2) BRUNETTI, F. Mecânica dos fluídos. São Paulo: Pearson Education. `n3) FOX, Robert W. Introdução à mecânica dos fluídos. Rio de Janeiro: LTC. `n4) WHITE, Frank M. Mecânica dos fluídos. Rio de Janeiro: Mcgraw-hill Interamericana. `n5) COULSON, J. M.; RICHARDSON, J.F. Chemical engineering . Oxford: Butterworth Heinemann. Volume 1: Fluid Flow, Heat Transfer and Mass Transfer"

# Narrow the first column group (was min=1,max=2 sharing style 1/width 30.71) so that
# column B gets its own distinct width (60.71) matching the target layout; this forces
# column A to be split off into its own <col> entry while keeping its original width.
$ws.Columns.Item(2).ColumnWidth = 59.83
